# Fix a typo in the PLL recognition guide: the paragraph that describes the
# "Gc" case ("-It's G?) ...") currently (incorrectly) reads "Gb" instead of
# "Gc". Only the FIRST occurrence of the phrase "Gb when the block is on the
# left" in the document (the one that sits under the "Gc:" heading) should be
# corrected - the other two occurrences (which belong to the "Gd:" heading,
# describing the Ga/Gc cases there) must stay untouched.

$d = $word.ActiveDocument

$needle = "Gb when the block is on the left"

$target = $null
$matchIndex = 0
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*" + $needle + "*") {
        $matchIndex = $matchIndex + 1
        if ($matchIndex -eq 1) {
            $target = $p.Range
        }
    }
}

if ($target -eq $null) {
    Write-Host "Target paragraph not found"
} else {
    $paraStart = $target.Start
    $paraText = $target.Text
    $localIdx = $paraText.IndexOf($needle)
    # Position of the "b" in "...s Gb when..." (right after "G")
    $bPos = $paraStart + $localIdx + 1
    $paraEnd = $target.End

    # --- Replace "b when the block is on the left" with
    #     "c when the block is on the left", going through a transient
    #     formatting change so the engine does not silently re-merge this
    #     text back into the preceding run (mirrors the run split Word
    #     itself produced for this edit).
    $rWhole = $d.Range($bPos, $paraEnd)
    $rWhole.Font.Bold = $true
    $rWhole.Text = "c when the block is on the left"
    $newEnd = $bPos + ("c when the block is on the left").Length
    $rWhole2 = $d.Range($bPos, $newEnd)
    $rWhole2.Font.Bold = $false

    # --- Now split "c" from " when the block is on the left" into two
    #     separate runs, again via a transient formatting toggle so the
    #     (now identically-formatted) pieces stay as distinct runs.
    $afterC = $bPos + 1
    $rSuffix = $d.Range($afterC, $newEnd)
    $rSuffix.Font.Bold = $true
    $rSuffix.Font.Bold = $false

    Write-Host "Replaced: $($d.Range($bPos - 1, $newEnd).Text)"
}
